$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the COVID-19 data row (row 6) with final figures
$ws.Range("C6").Value = 75670000
$ws.Range("D6").Value = 1670000
$ws.Range("E6").Value = 0.038
$ws.Range("F6").Value = 191
$ws.Range("G6").Value = (Get-Date -Year 2020 -Month 12 -Day 20 -Hour 0 -Minute 0 -Second 0).Date

# Update the view: scroll back to top and change selection
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 1
$ws.Range("F7").Select()
